$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename county name to statewide "Minnesota" for all data rows (2-19)
$ws.Range("A2:A19").Value = "Minnesota"

# Updated population figures (Total, Male, Female) per age band, rows 2-19
$data = @(
    @(337504, 172430, 165074),
    @(366552, 187432, 179120),
    @(382463, 196035, 186428),
    @(374886, 191061, 183825),
    @(353691, 178457, 175234),
    @(367229, 187130, 180099),
    @(381088, 195109, 185979),
    @(398280, 203093, 195187),
    @(368878, 189190, 179688),
    @(319464, 163330, 156134),
    @(345557, 175868, 169689),
    @(375195, 188800, 186395),
    @(380920, 189946, 190974),
    @(321627, 158506, 163121),
    @(255560, 122888, 132672),
    @(161437, 74725, 86712),
    @(107156, 46674, 60482),
    @(109903, 39739, 70164)
)

$row = 2
foreach ($vals in $data) {
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
    $row++
}
